# Auto-generated edit script applying numeric corrections to Sheets per commit diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2166.6667
$ws.Range("I28").Value = 2166.6667
$ws.Range("K28").Value = 2166.6667
$ws.Range("M28").Value = -1681.6667
$ws.Range("H51").Value = 5564.1333
$ws.Range("I51").Value = 5246.4443
$ws.Range("J51").Value = 6040.6665
$ws.Range("K51").Value = 5246.4443
$ws.Range("L51").Value = 6040.6665
$ws.Range("M51").Value = -4762.4443
$ws.Range("N51").Value = -7008.6665
$ws.Range("H103").Value = 62501036
$ws.Range("I103").Value = 125000980
$ws.Range("K103").Value = 375002940
$ws.Range("M103").Value = -375002354
$ws.Range("H111").Value = 2999.6667
$ws.Range("I111").Value = 3000
$ws.Range("K111").Value = 9000
$ws.Range("M111").Value = -5933
$ws.Range("H115").Value = 5998.8
$ws.Range("I115").Value = 5998.8
$ws.Range("K115").Value = 17996.4
$ws.Range("M115").Value = -16429.4
$ws.Range("H116").Value = 6300.8
$ws.Range("I116").Value = 5874.5
$ws.Range("J116").Value = 8006
$ws.Range("K116").Value = 5874.5
$ws.Range("L116").Value = 8006
$ws.Range("M116").Value = -2432.5
$ws.Range("N116").Value = -14890
$ws.Range("H135").Value = 3872.5
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").Value = $null
$ws.Range("H137").Value = 6231.647
$ws.Range("J137").Value = 6262.8667
$ws.Range("L137").Value = 18788.6001
$ws.Range("N137").Value = -23888.6001
$ws.Range("H138").Value = 2349.25
$ws.Range("I138").Value = 2266
$ws.Range("J138").Value = 2599
$ws.Range("K138").Value = 6798
$ws.Range("L138").Value = 7797
$ws.Range("M138").Value = -1658
$ws.Range("N138").Value = -18077

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 938.8043
$ws.Range("I32").Value = 968.093
$ws.Range("K32").Value = 968.093
$ws.Range("M32").Value = -681.093
$ws.Range("H45").Value = 1593
$ws.Range("I45").Value = 1476.7142
$ws.Range("J45").Value = 2000
$ws.Range("K45").Value = 1476.7142
$ws.Range("L45").Value = 2000
$ws.Range("M45").Value = -1099.7142
$ws.Range("N45").Value = -2754
$ws.Range("H49").Value = 25000
$ws.Range("J49").Value = 25000
$ws.Range("L49").Value = 25000
$ws.Range("N49").Value = -25520
$ws.Range("H61").Value = 2670.8484
$ws.Range("I61").Value = 2854.3447
$ws.Range("K61").Value = 2854.3447
$ws.Range("M61").Value = -2642.3447
$ws.Range("H63").Value = 3156.8604
$ws.Range("I63").Value = 1897.3287
$ws.Range("J63").Value = 10229.615
$ws.Range("K63").Value = 1897.3287
$ws.Range("L63").Value = 10229.615
$ws.Range("M63").Value = -1211.3287
$ws.Range("N63").Value = -11601.615
$ws.Range("H66").Value = 3156.8604
$ws.Range("I66").Value = 1897.3287
$ws.Range("J66").Value = 10229.615
$ws.Range("K66").Value = 9486.6435
$ws.Range("L66").Value = 51148.075
$ws.Range("M66").Value = -6054.6435
$ws.Range("N66").Value = -58012.075
$ws.Range("H132").Value = 1768.6
$ws.Range("I132").Value = 1760.2354
$ws.Range("K132").Value = 5280.706200000001
$ws.Range("M132").Value = -2750.706200000001
$ws.Range("H136").Value = 2670.8484
$ws.Range("I136").Value = 2854.3447
$ws.Range("K136").Value = 8563.034100000001
$ws.Range("M136").Value = -6013.034100000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1038.4546
$ws.Range("I20").Value = 946.1429000000001
$ws.Range("K20").Value = 946.1429000000001
$ws.Range("M20").Value = -699.1429000000001
$ws.Range("H94").Value = 384.5
$ws.Range("I94").Value = 384.5
$ws.Range("K94").Value = 384.5
$ws.Range("M94").Value = 66.5
$ws.Range("H99").Value = 27779990
$ws.Range("I99").Value = 11112351
$ws.Range("J99").Value = 55559390
$ws.Range("K99").Value = 11112351
$ws.Range("L99").Value = 55559390
$ws.Range("M99").Value = -11110853
$ws.Range("N99").Value = -55562386
$ws.Range("H107").Value = 1450.2354
$ws.Range("I107").Value = 1434.625
$ws.Range("K107").Value = 1434.625
$ws.Range("M107").Value = 485.375

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 3685
$ws.Range("I7").Value = 527.5
$ws.Range("J7").Value = 10000
$ws.Range("K7").Value = 527.5
$ws.Range("L7").Value = 10000
$ws.Range("M7").Value = -414.5
$ws.Range("N7").Value = -10226
$ws.Range("H58").Value = 2438.2632
$ws.Range("I58").Value = 1956.4286
$ws.Range("J58").Value = 3787.4
$ws.Range("K58").Value = 1956.4286
$ws.Range("L58").Value = 3787.4
$ws.Range("M58").Value = -1753.4286
$ws.Range("N58").Value = -4193.4
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").Value = $null
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").Value = $null
$ws.Range("H93").Value = 7302.3335
$ws.Range("I93").Value = 7302.3335
$ws.Range("K93").Value = 7302.3335
$ws.Range("M93").Value = -5430.3335
$ws.Range("H105").Value = 1050
$ws.Range("I105").Value = 900
$ws.Range("J105").Value = 1200
$ws.Range("K105").Value = 900
$ws.Range("L105").Value = 1200
$ws.Range("M105").Value = 847
$ws.Range("N105").Value = -4694
$ws.Range("H132").Value = 1919.4584
$ws.Range("I132").Value = 1846.2
$ws.Range("K132").Value = 5538.6
$ws.Range("M132").Value = -3008.6
$ws.Range("H134").Value = 2816.6191
$ws.Range("I134").Value = 2807.45
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 8422.349999999999
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -5887.349999999999
$ws.Range("N134").Value = -14070
$ws.Range("H136").Value = 2438.2632
$ws.Range("I136").Value = 1956.4286
$ws.Range("J136").Value = 3787.4
$ws.Range("K136").Value = 5869.2858
$ws.Range("L136").Value = 11362.2
$ws.Range("M136").Value = -3319.2858
$ws.Range("N136").Value = -16462.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1767.1111
$ws.Range("I11").Value = 816.6667
$ws.Range("J11").Value = 3668
$ws.Range("K11").Value = 2450.0001
$ws.Range("L11").Value = 11004
$ws.Range("M11").Value = -2310.0001
$ws.Range("N11").Value = -11284
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").Value = $null
$ws.Range("H70").Value = 10759.8
$ws.Range("I70").Value = 9033
$ws.Range("J70").Value = 13350
$ws.Range("K70").Value = 27099
$ws.Range("L70").Value = 40050
$ws.Range("M70").Value = -26784
$ws.Range("N70").Value = -40680
$ws.Range("H73").Value = 10759.8
$ws.Range("I73").Value = 9033
$ws.Range("J73").Value = 13350
$ws.Range("K73").Value = 27099
$ws.Range("L73").Value = 40050
$ws.Range("M73").Value = -26007
$ws.Range("N73").Value = -42234

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 379.85715
$ws.Range("I97").Value = 351.18182
$ws.Range("K97").Value = 351.18182
$ws.Range("M97").Value = 144.81818
$ws.Range("H102").Value = 1493.5
$ws.Range("I102").Value = 1591.4
$ws.Range("K102").Value = 1591.4
$ws.Range("M102").Value = 30.59999999999991
$ws.Range("H113").Value = 1426
$ws.Range("I113").Value = 1426
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1426
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 744
$ws.Range("N113").Value = $null
$ws.Range("H136").Value = 40326
$ws.Range("J136").Value = 40326
$ws.Range("L136").Value = 120978
$ws.Range("N136").Value = -126078

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9099.200000000001
$ws.Range("I40").Value = 4554.6665
$ws.Range("J40").Value = 50000
$ws.Range("K40").Value = 4554.6665
$ws.Range("L40").Value = 50000
$ws.Range("M40").Value = -4418.6665
$ws.Range("N40").Value = -50272
$ws.Range("H105").Value = 22538.334
$ws.Range("J105").Value = 22538.334
$ws.Range("L105").Value = 22538.334
$ws.Range("N105").Value = -29526.334
$ws.Range("H122").Value = 18374.25
$ws.Range("I122").Value = 17400
$ws.Range("K122").Value = 52200
$ws.Range("M122").Value = -49750

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 59500
$ws.Range("J112").Value = 59500
$ws.Range("L112").Value = 59500
$ws.Range("N112").Value = -62454
$ws.Range("H122").Value = 2414.5715
$ws.Range("I122").Value = 1806.2354
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 5418.706200000001
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -2968.706200000001
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 3898.1667
$ws.Range("I132").Value = 2743.3076
$ws.Range("K132").Value = 8229.9228
$ws.Range("M132").Value = -5699.9228
$ws.Range("H136").Value = 1902.6154
$ws.Range("I136").Value = 1294
$ws.Range("K136").Value = 3882
$ws.Range("M136").Value = -1332
